$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.984.47'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '2.264.15'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.650'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.71'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.68%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.450'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0985'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").Value = '2.601.59'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("E16").Value = '  +3.72%  '
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("D18").Value = '2.268.92'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").Value = '43.932.64'
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("E20").Value = '  +4.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  -3.64%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.66%  '
$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +21.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.19%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("E33").Value = '  +2.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0688'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.70'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.48'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.85%  '
$ws.Range("E39").Value = '  -2.43%  '
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000223'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.92%  '
$ws.Range("E44").Value = '  +2.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '98.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("E47").Value = '  -1.30%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.79%  '
$ws.Range("D50").Value = '1.448.82'
$ws.Range("E50").Value = '  -2.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.11%  '
